$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new list paragraph "MSGEX13 Datos incorrectos." right
#    after the "MSGEX12 Receta no creada." paragraph (same list /
#    numbering as its neighbours), mirroring the bold "MSGEXxx" +
#    normal description run pattern used throughout that list.
# ------------------------------------------------------------------
$msgex12 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*MSGEX12 Receta no creada.*") {
        $msgex12 = $p
    }
}

if ($msgex12 -ne $null) {
    # Grab the formatting (bold + complex-script bold) of the existing
    # "MSGEX12" label run so the new label picks up identical run
    # properties (<w:b/><w:bCs/> etc.) instead of only <w:b/>.
    $labelSrc = $d.Range($msgex12.Range.Start, $msgex12.Range.Start + 7)
    $labelFormatted = $labelSrc.FormattedText

    # New paragraph inherits pPr/numPr/style from msgex12 automatically.
    $msgex12.Range.InsertParagraphAfter()
    $newPara = $msgex12.Next()
    $newRange = $newPara.Range
    $insertStart = $newRange.Start

    # Seed the final text first (plain run), then re-apply the bold
    # label formatting over the label's character span.
    $newRange.Text = "MSGEX13 Datos incorrectos."

    $labelDest = $d.Range($insertStart, $insertStart + 7)
    $labelDest.FormattedText = $labelFormatted

    # FormattedText carried the *source* characters ("MSGEX12") along
    # with its formatting; restore the correct label text in place
    # (same run/formatting, only the literal characters change).
    $fixRange = $d.Range($insertStart + 6, $insertStart + 7)
    $fixRange.Text = "3"
}

# ------------------------------------------------------------------
# 2) "Medicamento agregado correctamente." -> "Medicamento agregado."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Medicamento agregado correctamente.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Medicamento agregado.", 2) | Out-Null
